$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G (pushes the old F..M data over to H..O)
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("F1").Value = "reviewed_at"
$ws.Range("G1").Value = "escalated_to"

# Populate the new columns for each data row
$ws.Range("F2").Value = "Division Performance"
$ws.Range("G2").Value = "Trust Performance"

$ws.Range("F3").Value = "Division Performance"
$ws.Range("G3").Value = "Trust Performance"

$ws.Range("F4").Value = "Division Performance"
$ws.Range("G4").Value = "Trust Performance"

$ws.Range("F5").Value = "Division Performance"
$ws.Range("G5").Value = "Trust Performance"

$ws.Range("F6").Value = "HR Review"
$ws.Range("G6").Value = "Board"

$ws.Range("F7").Value = "Directorate Performance"
$ws.Range("G7").Value = "Division Performance"

# Autofit the new columns to their content, like the other data columns
$ws.Range("F1:G7").EntireColumn.AutoFit() | Out-Null

# Restore single-cell selection (matches the saved view state)
$ws.Range("G14").Select() | Out-Null
